$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells (rows 2, 4, 9) ---
# Row 2
$ws.Range("M2").Value = 36
$ws.Range("O2").Value = 922.1

# Row 4
$ws.Range("K4").Value = 35
$ws.Range("M4").Value = 310
$ws.Range("O4").Value = 2324.7

# Row 9
$ws.Range("K9").Value = 87
$ws.Range("M9").Value = 1657
$ws.Range("O9").Value = 2604.2

# --- Add new rows 10-13, matching the layout/styles of row 9 but without column E ---
$ws.Range("A9:O9").Copy()
$ws.Range("A10:O13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Clear the column E cells on the new rows (they are not used there)
$ws.Range("E10:E13").Clear()

# Row 10
$ws.Range("A10").Value = "local"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = "http_test_case_1"
$ws.Range("F10").Value = 20
$ws.Range("G10").Value = 10
$ws.Range("H10").Value = 1000
$ws.Range("I10").Value = 10
$ws.Range("J10").Value = "/callbackNoDapr"
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 19
$ws.Range("O10").Value = 973

# Row 11
$ws.Range("A11").Value = "local"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "http_test_case_1"
$ws.Range("F11").Value = 20
$ws.Range("G11").Value = 10
$ws.Range("H11").Value = 1000
$ws.Range("I11").Value = 50
$ws.Range("J11").Value = "/callbackNoDapr"
$ws.Range("K11").Value = 14
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 127
$ws.Range("O11").Value = 2280.9

# Row 12
$ws.Range("A12").Value = "local"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = "http_test_case_1"
$ws.Range("F12").Value = 20
$ws.Range("G12").Value = 10
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 100
$ws.Range("J12").Value = "/callbackNoDapr"
$ws.Range("K12").Value = 25
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 159
$ws.Range("O12").Value = 2924.3

# Row 13
$ws.Range("A13").Value = "local"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = "http_test_case_1"
$ws.Range("F13").Value = 20
$ws.Range("G13").Value = 10
$ws.Range("H13").Value = 1000
$ws.Range("I13").Value = 250
$ws.Range("J13").Value = "/callbackNoDapr"
$ws.Range("K13").Value = 84
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 517
$ws.Range("O13").Value = 2722

# --- Resize the table and its autofilter to include the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:O13"))

# --- Column F width: 16.5 characters (customWidth) ---
$ws.Columns.Item(6).ColumnWidth = 15.666666666666666

# --- Update the active selection ---
$ws.Range("J10").Select() | Out-Null
